$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing E column values (E2:E5)
$ws.Range("E2").Value = 3197.400331497192
$ws.Range("E3").Value = 2525.983095169067
$ws.Range("E4").Value = 2914.046764373779
$ws.Range("E5").Value = 2470.466375350952

# Append new rows (6-11) with data for generations 4-9
$newRows = @(
    @(4, "prey wins", 100, 20, 2393.391847610474),
    @(5, "prey wins", 100, 20, 2444.068670272827),
    @(6, "prey wins", 100, 20, 2441.547632217407),
    @(7, "prey wins", 100, 20, 2436.452627182007),
    @(8, "prey wins", 100, 20, 2453.410387039185),
    @(9, "prey wins", 100, 20, 2459.134340286255)
)

$rowIndex = 6
foreach ($rowData in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $rowIndex++
}
